$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("December")

# New entries for December (row 4 through row 8)
$ws.Range("C4").Value = "Mas+ Peaj +roshun+riska"
$ws.Range("D4").Value = 750
$ws.Range("E4").Value = "shim "
$ws.Range("F4").Value = 35

$ws.Range("E5").Value = "kopi"
$ws.Range("F5").Value = 30

$ws.Range("E6").Value = "peaj pata "
$ws.Range("F6").Value = 25

$ws.Range("E7").Value = "morich"
$ws.Range("F7").Value = 10

$ws.Range("E8").Value = "Mas "
$ws.Range("F8").Value = 110

# Move the active selection to F9, matching the saved view state
$ws.Range("F9").Select()
